$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NCAP_BND")

# Rename the UD set reference in the "Attribute" header row from PSET_CI to Pset_set
$ws.Range("D18").Value = "Pset_set"

# Rename ELCGAS -> ELEGAS across the data rows (19-25)
$ws.Range("D19:D25").Value = "ELEGAS"

# Update the active selection to D20 (as last saved by the author)
$ws.Range("D20").Select()
